$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 4).Value = 44204
$ws.Cells.Item(3, 11).Value = 'Bing'
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 400
$ws.Cells.Item(3, 14).Value = 10000
$ws.Cells.Item(3, 15).Value = 11000
$ws.Cells.Item(3, 16).Value = 10500
$ws.Cells.Item(3, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(3, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(3, 19).Value = 1050
$ws.Cells.Item(3, 20).Value = 10

# Row 4
$ws.Cells.Item(4, 4).Value = 44204
$ws.Cells.Item(4, 11).Value = 'Rainier'
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 400
$ws.Cells.Item(4, 14).Value = 15000
$ws.Cells.Item(4, 15).Value = 16000
$ws.Cells.Item(4, 16).Value = 15500
$ws.Cells.Item(4, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(4, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(4, 19).Value = 1550
$ws.Cells.Item(4, 20).Value = 10

# Row 5
$ws.Cells.Item(5, 4).Value = 44211
$ws.Cells.Item(5, 11).Value = 'Bing'
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 400
$ws.Cells.Item(5, 14).Value = 6000
$ws.Cells.Item(5, 15).Value = 7000
$ws.Cells.Item(5, 16).Value = 6500
$ws.Cells.Item(5, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(5, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(5, 19).Value = 929
$ws.Cells.Item(5, 20).Value = 7

# Row 6
$ws.Cells.Item(6, 4).Value = 44211
$ws.Cells.Item(6, 11).Value = 'Rainier'
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 400
$ws.Cells.Item(6, 14).Value = 9000
$ws.Cells.Item(6, 15).Value = 10000
$ws.Cells.Item(6, 16).Value = 9500
$ws.Cells.Item(6, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(6, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(6, 19).Value = 1357
$ws.Cells.Item(6, 20).Value = 7

# Row 7
$ws.Cells.Item(7, 4).Value = 44159
$ws.Cells.Item(7, 11).Value = 'Early Burlat'
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 500
$ws.Cells.Item(7, 14).Value = 11000
$ws.Cells.Item(7, 15).Value = 12000
$ws.Cells.Item(7, 16).Value = 11500
$ws.Cells.Item(7, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(7, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7, 19).Value = 1643
$ws.Cells.Item(7, 20).Value = 7

# Row 8
$ws.Cells.Item(8, 4).Value = 44523
$ws.Cells.Item(8, 11).Value = 'Santina'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 400
$ws.Cells.Item(8, 14).Value = 27000
$ws.Cells.Item(8, 15).Value = 28000
$ws.Cells.Item(8, 16).Value = 27500
$ws.Cells.Item(8, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(8, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(8, 19).Value = 2750
$ws.Cells.Item(8, 20).Value = 10

# Row 9
$ws.Cells.Item(9, 4).Value = 44186
$ws.Cells.Item(9, 11).Value = 'Rainier'
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 400
$ws.Cells.Item(9, 14).Value = 15000
$ws.Cells.Item(9, 15).Value = 16000
$ws.Cells.Item(9, 16).Value = 15500
$ws.Cells.Item(9, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(9, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(9, 19).Value = 1550
$ws.Cells.Item(9, 20).Value = 10

# Row 10
$ws.Cells.Item(10, 4).Value = 44189
$ws.Cells.Item(10, 11).Value = 'Bing'
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 500
$ws.Cells.Item(10, 14).Value = 11000
$ws.Cells.Item(10, 15).Value = 12000
$ws.Cells.Item(10, 16).Value = 11500
$ws.Cells.Item(10, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(10, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(10, 19).Value = 1150
$ws.Cells.Item(10, 20).Value = 10

# Row 11
$ws.Cells.Item(11, 4).Value = 44189
$ws.Cells.Item(11, 11).Value = 'Rainier'
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 500
$ws.Cells.Item(11, 14).Value = 16000
$ws.Cells.Item(11, 15).Value = 17000
$ws.Cells.Item(11, 16).Value = 16500
$ws.Cells.Item(11, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(11, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(11, 19).Value = 1650
$ws.Cells.Item(11, 20).Value = 10

# Row 12
$ws.Cells.Item(12, 4).Value = 44169
$ws.Cells.Item(12, 11).Value = 'Lapins'
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 600
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 15500
$ws.Cells.Item(12, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 1550
$ws.Cells.Item(12, 20).Value = 10

# Row 13
$ws.Cells.Item(13, 4).Value = 44215
$ws.Cells.Item(13, 11).Value = 'Bing'
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 300
$ws.Cells.Item(13, 14).Value = 6000
$ws.Cells.Item(13, 15).Value = 7000
$ws.Cells.Item(13, 16).Value = 6500
$ws.Cells.Item(13, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(13, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(13, 19).Value = 929
$ws.Cells.Item(13, 20).Value = 7

# Row 14
$ws.Cells.Item(14, 4).Value = 44215
$ws.Cells.Item(14, 11).Value = 'Rainier'
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 9000
$ws.Cells.Item(14, 15).Value = 10000
$ws.Cells.Item(14, 16).Value = 9500
$ws.Cells.Item(14, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(14, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(14, 19).Value = 1357
$ws.Cells.Item(14, 20).Value = 7

# Row 15
$ws.Cells.Item(15, 4).Value = 44201
$ws.Cells.Item(15, 11).Value = 'Bing'
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 600
$ws.Cells.Item(15, 14).Value = 11000
$ws.Cells.Item(15, 15).Value = 12000
$ws.Cells.Item(15, 16).Value = 11500
$ws.Cells.Item(15, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(15, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(15, 19).Value = 1150
$ws.Cells.Item(15, 20).Value = 10

# Row 16
$ws.Cells.Item(16, 4).Value = 44201
$ws.Cells.Item(16, 11).Value = 'Rainier'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 600
$ws.Cells.Item(16, 14).Value = 16000
$ws.Cells.Item(16, 15).Value = 17000
$ws.Cells.Item(16, 16).Value = 16500
$ws.Cells.Item(16, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(16, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(16, 19).Value = 1650
$ws.Cells.Item(16, 20).Value = 10

# Row 17
$ws.Cells.Item(17, 4).Value = 44530
$ws.Cells.Item(17, 11).Value = 'Santina'
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 600
$ws.Cells.Item(17, 14).Value = 20000
$ws.Cells.Item(17, 15).Value = 21000
$ws.Cells.Item(17, 16).Value = 20500
$ws.Cells.Item(17, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(17, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(17, 19).Value = 1367
$ws.Cells.Item(17, 20).Value = 15

# Row 18
$ws.Cells.Item(18, 4).Value = 44537
$ws.Cells.Item(18, 11).Value = 'Santina'
$ws.Cells.Item(18, 12).Value = 'Primera'
$ws.Cells.Item(18, 13).Value = 600
$ws.Cells.Item(18, 14).Value = 10000
$ws.Cells.Item(18, 15).Value = 11000
$ws.Cells.Item(18, 16).Value = 10500
$ws.Cells.Item(18, 17).Value = '$/caja 8 kilos'
$ws.Cells.Item(18, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(18, 19).Value = 1312
$ws.Cells.Item(18, 20).Value = 8

# Row 19
$ws.Cells.Item(19, 4).Value = 44162
$ws.Cells.Item(19, 11).Value = 'Early Burlat'
$ws.Cells.Item(19, 12).Value = 'Segunda'
$ws.Cells.Item(19, 13).Value = 500
$ws.Cells.Item(19, 14).Value = 11000
$ws.Cells.Item(19, 15).Value = 12000
$ws.Cells.Item(19, 16).Value = 11500
$ws.Cells.Item(19, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(19, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 19).Value = 1643
$ws.Cells.Item(19, 20).Value = 7

# Row 20
$ws.Cells.Item(20, 4).Value = 44222
$ws.Cells.Item(20, 11).Value = 'Lapins'
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 200
$ws.Cells.Item(20, 14).Value = 6000
$ws.Cells.Item(20, 15).Value = 7000
$ws.Cells.Item(20, 16).Value = 6500
$ws.Cells.Item(20, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(20, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(20, 19).Value = 929
$ws.Cells.Item(20, 20).Value = 7

# Row 21
$ws.Cells.Item(21, 4).Value = 44166
$ws.Cells.Item(21, 11).Value = 'Lapins'
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 700
$ws.Cells.Item(21, 14).Value = 15000
$ws.Cells.Item(21, 15).Value = 16000
$ws.Cells.Item(21, 16).Value = 15500
$ws.Cells.Item(21, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(21, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(21, 19).Value = 1550
$ws.Cells.Item(21, 20).Value = 10

# Row 22
$ws.Cells.Item(22, 4).Value = 44516
$ws.Cells.Item(22, 11).Value = 'Early Burlat'
$ws.Cells.Item(22, 12).Value = 'Segunda'
$ws.Cells.Item(22, 13).Value = 600
$ws.Cells.Item(22, 14).Value = 13000
$ws.Cells.Item(22, 15).Value = 13500
$ws.Cells.Item(22, 16).Value = 13250
$ws.Cells.Item(22, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(22, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(22, 19).Value = 2650
$ws.Cells.Item(22, 20).Value = 5

# Row 23
$ws.Cells.Item(23, 4).Value = 44533
$ws.Cells.Item(23, 11).Value = 'Santina'
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 600
$ws.Cells.Item(23, 14).Value = 10000
$ws.Cells.Item(23, 15).Value = 11000
$ws.Cells.Item(23, 16).Value = 10500
$ws.Cells.Item(23, 17).Value = '$/caja 8 kilos'
$ws.Cells.Item(23, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(23, 19).Value = 1312
$ws.Cells.Item(23, 20).Value = 8

# Row 24
$ws.Cells.Item(24, 4).Value = 44176
$ws.Cells.Item(24, 11).Value = 'Bing'
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 500
$ws.Cells.Item(24, 14).Value = 7000
$ws.Cells.Item(24, 15).Value = 8000
$ws.Cells.Item(24, 16).Value = 7500
$ws.Cells.Item(24, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(24, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(24, 19).Value = 1500
$ws.Cells.Item(24, 20).Value = 5

# Row 25
$ws.Cells.Item(25, 4).Value = 44176
$ws.Cells.Item(25, 11).Value = 'Lapins'
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 500
$ws.Cells.Item(25, 14).Value = 13000
$ws.Cells.Item(25, 15).Value = 14000
$ws.Cells.Item(25, 16).Value = 13500
$ws.Cells.Item(25, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(25, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(25, 19).Value = 1350
$ws.Cells.Item(25, 20).Value = 10

# Row 26
$ws.Cells.Item(26, 4).Value = 44539
$ws.Cells.Item(26, 11).Value = 'Santina'
$ws.Cells.Item(26, 12).Value = 'Primera'
$ws.Cells.Item(26, 13).Value = 200
$ws.Cells.Item(26, 14).Value = 11000
$ws.Cells.Item(26, 15).Value = 12000
$ws.Cells.Item(26, 16).Value = 11500
$ws.Cells.Item(26, 17).Value = '$/caja 8 kilos'
$ws.Cells.Item(26, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(26, 19).Value = 1438
$ws.Cells.Item(26, 20).Value = 8

# Row 27
$ws.Cells.Item(27, 4).Value = 44532
$ws.Cells.Item(27, 11).Value = 'Santina'
$ws.Cells.Item(27, 12).Value = 'Primera'
$ws.Cells.Item(27, 13).Value = 200
$ws.Cells.Item(27, 14).Value = 20000
$ws.Cells.Item(27, 15).Value = 21000
$ws.Cells.Item(27, 16).Value = 20500
$ws.Cells.Item(27, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(27, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(27, 19).Value = 1367
$ws.Cells.Item(27, 20).Value = 15

# Row 28
$ws.Cells.Item(28, 4).Value = 44218
$ws.Cells.Item(28, 11).Value = 'Lapins'
$ws.Cells.Item(28, 12).Value = 'Primera'
$ws.Cells.Item(28, 13).Value = 300
$ws.Cells.Item(28, 14).Value = 6000
$ws.Cells.Item(28, 15).Value = 7000
$ws.Cells.Item(28, 16).Value = 6500
$ws.Cells.Item(28, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(28, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(28, 19).Value = 929
$ws.Cells.Item(28, 20).Value = 7

# Row 29
$ws.Cells.Item(29, 4).Value = 44161
$ws.Cells.Item(29, 11).Value = 'Early Burlat'
$ws.Cells.Item(29, 12).Value = 'Segunda'
$ws.Cells.Item(29, 13).Value = 400
$ws.Cells.Item(29, 14).Value = 11000
$ws.Cells.Item(29, 15).Value = 12000
$ws.Cells.Item(29, 16).Value = 11500
$ws.Cells.Item(29, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(29, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(29, 19).Value = 1643
$ws.Cells.Item(29, 20).Value = 7

# Row 30
$ws.Cells.Item(30, 4).Value = 44208
$ws.Cells.Item(30, 11).Value = 'Bing'
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 300
$ws.Cells.Item(30, 14).Value = 6000
$ws.Cells.Item(30, 15).Value = 7000
$ws.Cells.Item(30, 16).Value = 6500
$ws.Cells.Item(30, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(30, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(30, 19).Value = 929
$ws.Cells.Item(30, 20).Value = 7

# Row 31
$ws.Cells.Item(31, 4).Value = 44208
$ws.Cells.Item(31, 11).Value = 'Rainier'
$ws.Cells.Item(31, 12).Value = 'Primera'
$ws.Cells.Item(31, 13).Value = 300
$ws.Cells.Item(31, 14).Value = 9000
$ws.Cells.Item(31, 15).Value = 10000
$ws.Cells.Item(31, 16).Value = 9500
$ws.Cells.Item(31, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(31, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(31, 19).Value = 1357
$ws.Cells.Item(31, 20).Value = 7

# Row 32
$ws.Cells.Item(32, 4).Value = 44519
$ws.Cells.Item(32, 11).Value = 'Early Burlat'
$ws.Cells.Item(32, 12).Value = 'Segunda'
$ws.Cells.Item(32, 13).Value = 600
$ws.Cells.Item(32, 14).Value = 13000
$ws.Cells.Item(32, 15).Value = 13500
$ws.Cells.Item(32, 16).Value = 13250
$ws.Cells.Item(32, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(32, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(32, 19).Value = 2650
$ws.Cells.Item(32, 20).Value = 5

# Row 33
$ws.Cells.Item(33, 4).Value = 44187
$ws.Cells.Item(33, 11).Value = 'Lapins'
$ws.Cells.Item(33, 12).Value = 'Primera'
$ws.Cells.Item(33, 13).Value = 600
$ws.Cells.Item(33, 14).Value = 9000
$ws.Cells.Item(33, 15).Value = 10000
$ws.Cells.Item(33, 16).Value = 9500
$ws.Cells.Item(33, 17).Value = '$/caja 8 kilos'
$ws.Cells.Item(33, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(33, 19).Value = 1188
$ws.Cells.Item(33, 20).Value = 8

# Row 34
$ws.Cells.Item(34, 4).Value = 44187
$ws.Cells.Item(34, 11).Value = 'Rainier'
$ws.Cells.Item(34, 12).Value = 'Primera'
$ws.Cells.Item(34, 13).Value = 600
$ws.Cells.Item(34, 14).Value = 16000
$ws.Cells.Item(34, 15).Value = 17000
$ws.Cells.Item(34, 16).Value = 16500
$ws.Cells.Item(34, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(34, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(34, 19).Value = 1650
$ws.Cells.Item(34, 20).Value = 10

# Row 35
$ws.Cells.Item(35, 4).Value = 44196
$ws.Cells.Item(35, 11).Value = 'Bing'
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 12000
$ws.Cells.Item(35, 15).Value = 13000
$ws.Cells.Item(35, 16).Value = 12500
$ws.Cells.Item(35, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(35, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(35, 19).Value = 1250
$ws.Cells.Item(35, 20).Value = 10

# Row 36
$ws.Cells.Item(36, 4).Value = 44196
$ws.Cells.Item(36, 11).Value = 'Rainier'
$ws.Cells.Item(36, 12).Value = 'Primera'
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 16000
$ws.Cells.Item(36, 15).Value = 17000
$ws.Cells.Item(36, 16).Value = 16500
$ws.Cells.Item(36, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(36, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(36, 19).Value = 1650
$ws.Cells.Item(36, 20).Value = 10

# Row 37
$ws.Cells.Item(37, 4).Value = 44200
$ws.Cells.Item(37, 11).Value = 'Bing'
$ws.Cells.Item(37, 12).Value = 'Primera'
$ws.Cells.Item(37, 13).Value = 200
$ws.Cells.Item(37, 14).Value = 12000
$ws.Cells.Item(37, 15).Value = 13000
$ws.Cells.Item(37, 16).Value = 12500
$ws.Cells.Item(37, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(37, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(37, 19).Value = 1250
$ws.Cells.Item(37, 20).Value = 10

# Row 38
$ws.Cells.Item(38, 4).Value = 44195
$ws.Cells.Item(38, 11).Value = 'Bing'
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 600
$ws.Cells.Item(38, 14).Value = 12000
$ws.Cells.Item(38, 15).Value = 13000
$ws.Cells.Item(38, 16).Value = 12500
$ws.Cells.Item(38, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(38, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(38, 19).Value = 1250
$ws.Cells.Item(38, 20).Value = 10

# Row 39
$ws.Cells.Item(39, 4).Value = 44195
$ws.Cells.Item(39, 11).Value = 'Rainier'
$ws.Cells.Item(39, 12).Value = 'Primera'
$ws.Cells.Item(39, 13).Value = 600
$ws.Cells.Item(39, 14).Value = 16000
$ws.Cells.Item(39, 15).Value = 17000
$ws.Cells.Item(39, 16).Value = 16500
$ws.Cells.Item(39, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(39, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(39, 19).Value = 1650
$ws.Cells.Item(39, 20).Value = 10

# Row 40
$ws.Cells.Item(40, 4).Value = 44536
$ws.Cells.Item(40, 11).Value = 'Santina'
$ws.Cells.Item(40, 12).Value = 'Primera'
$ws.Cells.Item(40, 13).Value = 400
$ws.Cells.Item(40, 14).Value = 10000
$ws.Cells.Item(40, 15).Value = 11000
$ws.Cells.Item(40, 16).Value = 10500
$ws.Cells.Item(40, 17).Value = '$/bandeja 8 kilos'
$ws.Cells.Item(40, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(40, 19).Value = 1312
$ws.Cells.Item(40, 20).Value = 8

# Row 41
$ws.Cells.Item(41, 4).Value = 44194
$ws.Cells.Item(41, 11).Value = 'Bing'
$ws.Cells.Item(41, 12).Value = 'Primera'
$ws.Cells.Item(41, 13).Value = 500
$ws.Cells.Item(41, 14).Value = 11000
$ws.Cells.Item(41, 15).Value = 12000
$ws.Cells.Item(41, 16).Value = 11500
$ws.Cells.Item(41, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(41, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(41, 19).Value = 1150
$ws.Cells.Item(41, 20).Value = 10

# Row 42
$ws.Cells.Item(42, 4).Value = 44194
$ws.Cells.Item(42, 11).Value = 'Rainier'
$ws.Cells.Item(42, 12).Value = 'Primera'
$ws.Cells.Item(42, 13).Value = 600
$ws.Cells.Item(42, 14).Value = 16000
$ws.Cells.Item(42, 15).Value = 17000
$ws.Cells.Item(42, 16).Value = 16500
$ws.Cells.Item(42, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(42, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(42, 19).Value = 1650
$ws.Cells.Item(42, 20).Value = 10

# Row 43
$ws.Cells.Item(43, 4).Value = 44518
$ws.Cells.Item(43, 11).Value = 'Early Burlat'
$ws.Cells.Item(43, 12).Value = 'Segunda'
$ws.Cells.Item(43, 13).Value = 200
$ws.Cells.Item(43, 14).Value = 13000
$ws.Cells.Item(43, 15).Value = 13500
$ws.Cells.Item(43, 16).Value = 13250
$ws.Cells.Item(43, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(43, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(43, 19).Value = 2650
$ws.Cells.Item(43, 20).Value = 5

# Row 44
$ws.Cells.Item(44, 4).Value = 44540
$ws.Cells.Item(44, 11).Value = 'Santina'
$ws.Cells.Item(44, 12).Value = 'Primera'
$ws.Cells.Item(44, 13).Value = 800
$ws.Cells.Item(44, 14).Value = 9000
$ws.Cells.Item(44, 15).Value = 9000
$ws.Cells.Item(44, 16).Value = 9000
$ws.Cells.Item(44, 17).Value = '$/caja 8 kilos'
$ws.Cells.Item(44, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(44, 19).Value = 1125
$ws.Cells.Item(44, 20).Value = 8

Write-Host "Done applying edits"